# "Getting Started" guide update
#
#  - The "Prepare a lesson plan" step's blurb about worksheets is reworded:
#      ", and is a complete project to demonstrate an aspect of machine
#      learning."
#    becomes
#      ", and is a complete project to demonstrate a real-world use of
#      machine learning."
#    and the sentence about downloading a template Scratch project file is
#    dropped (it is no longer needed, now templates can be downloaded from
#    inside Scratch) - this sentence is removed from BOTH places it occurs
#    in the document.
#  - The (hidden) "_GoBack" bookmark that Word leaves behind at the last
#    edit position moves from the "unmanaged class accounts" guide (where
#    it is simply deleted) to the "managed class accounts" guide, landing
#    right after the (now-unchanged) "...aspect of machine learning. "
#    sentence, where the removed sentence used to be.
#  - The footer's "Last updated" date is bumped.

$d = $word.ActiveDocument

# wdFindContinue=1, wdReplaceOne=1, wdReplaceAll=2, wdCollapseEnd=0 (used implicitly)

# --- Update the "Last updated" date in the footer ---------------------
foreach ($story in $d.StoryRanges) {
    if ($story.Text -like "*Last updated:*") {
        $story.Find.Execute("18 October 2018", $true, $false, $false, $false, `
            $false, $true, 1, $false, "11 November 2018", 2) | Out-Null
    }
}

# --- Drop the old (hidden) _GoBack bookmark ----------------------------
# It currently sits just before "You will need to do this if you want your
# class to be able to do projects that recognise images" in the "unmanaged
# class accounts" guide. It gets removed from here, and re-created later
# at its new home.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# --- First occurrence: reword the sentence, in the "managed class
#     accounts" guide's "Prepare a lesson plan" step ---------------------
$rWording = $d.Content
$rWording.Find.Execute( `
    ", and is a complete project to demonstrate an aspect of machine learning. ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    ", and is a complete project to demonstrate a real-world use of machine learning. ", `
    1) | Out-Null

# Remove the (now obsolete) sentence about template Scratch project files,
# along with the manual line break in front of it, for this first
# occurrence only.
$vtab = [char]11
$templateSentence = $vtab + "Some of the projects include a template Scratch project file – these are available for download alongside the worksheet PDFs. "

$rFirstCut = $d.Content
$rFirstCut.Find.Execute($templateSentence) | Out-Null
$rFirstCut.Delete()

# --- Second occurrence: in the "unmanaged class accounts" guide's
#     "Prepare a lesson plan" step, the wording of the main sentence is
#     unchanged, but the template-Scratch-file sentence (and its leading
#     line break) is removed and replaced with the relocated _GoBack
#     bookmark. ------------------------------------------------------
$rSecondCut = $d.Content
$rSecondCut.Find.Execute($templateSentence) | Out-Null
$rSecondCut.Delete()
$d.Bookmarks.Add("_GoBack", $rSecondCut) | Out-Null
